# Rename the "Test issueN" values in column A to "Test issueXLSXN"
# (data in columns B-E stays identical; only the labels in column A change)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = "Test issueXLSX$r"
}

# Widen column A to fit the longer labels (stored width ends up as 15)
$ws.Columns.Item(1).ColumnWidth = 14.166666666666666

# Move the active selection from G7 to C1
$ws.Range("C1").Select()
